{"js": "// Apply the LOQ4206 course-plan updates:\n//  1. Activation date 2018 -> 2021\n//  2. Second \"Programa\" English summary paragraph expanded with full syllabus detail\n//  3. \"M\u00e9todo\" evaluation description replaced\n//  4. \"Crit\u00e9rio\" formula replaced\n//  5. \"Norma de recupera\u00e7\u00e3o\" text replaced\n\nconst body = context.document.body;\n\n// 1) Ativa\u00e7\u00e3o date -------------------------------------------------------\nconst ativacao = body.search(\"Ativa\u00e7\u00e3o: 01/01/2018\", { matchCase: true });\nativacao.load(\"items\");\nawait context.sync();\nif (ativacao.items.length > 0) {\n  ativacao.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2021\", Word.InsertLocation.replace);\n}\n\n// 2) Expand the second occurrence of the English \"Programa\" summary -----\n// (identical text also appears once under \"Programa resumido\" - that one\n// must stay untouched, so we locate it via the paragraph collection and\n// update only the *second* match.)\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst shortProgramaText =\n  \"Introduction to Operational Research, Linear Programming, Simplex Method, \" +\n  \"Introduction to Graphs and Network Optimization, Case Study in Linear \" +\n  \"Programming, Introduction to Queue Theory.\";\nconst fullProgramaText =\n  \"1. Introduction to Operational Research 1.1. Concepts of Operational Research; \" +\n  \"1.2. Modeling; 1.3. Structure of Mathematical Models; 1.4. Mathematical techniques \" +\n  \"in Operational Research; 1.2. Phases of a Study in Operational Research 2. Linear \" +\n  \"Programming 2.1. Definition 2.2. Formulation of Models 2.3. Graphic Resolution; \" +\n  \"3. Simplex method 3.1. Development of the Simplex Method; 3.2. Simplex Method \" +\n  \"Procedure; 4. Introduction to Graphs and Network Optimization 4.1. Basic Concepts \" +\n  \"in Graph Theory 4.2. Maximum Flow Problems; 4.3. Minimum Path Problems 5. Case \" +\n  \"Studies in Linear Programming 5.1. Simple Transport Model 5.2. Model of \" +\n  \"Designation. 6. Introduction to Queuing Theory 6.1. Queuing Theory Concepts \" +\n  \"6.2. Markovian Models\";\n\nlet seen = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === shortProgramaText) {\n    seen++;\n    if (seen === 2) {\n      paragraphs.items[i].getRange().insertText(fullProgramaText, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\n\n// 3) M\u00e9todo description --------------------------------------------------\nconst metodo = body.search(\n  \"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de laborat\u00f3rio, aulas de exerc\u00edcios.\",\n  { matchCase: true }\n);\nmetodo.load(\"items\");\nawait context.sync();\nif (metodo.items.length > 0) {\n  metodo.items[0].insertText(\n    \"NF=A avalia\u00e7\u00e3o ser\u00e1 composta por provas, listas, projetos, semin\u00e1rios e outras formas que far\u00e3o a composi\u00e7\u00e3o das notas, sendo estipulada a m\u00e9dia final a somat\u00f3ria destas notas (N), com no m\u00ednimo duas avalia\u00e7\u00f5es, sendo: (N1+...+Nn)/n.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 4) Crit\u00e9rio formula ------------------------------------------------------\nconst criterio = body.search(\n  \"MF = (0,45*P1 + 0,45*P2 + 0,10*TRAB), onde P1 e P2 s\u00e3o provas e TRAB \u00e9 a nota m\u00e9dia de trabalhos e semin\u00e1rios.\",\n  { matchCase: true }\n);\ncriterio.load(\"items\");\nawait context.sync();\nif (criterio.items.length > 0) {\n  criterio.items[0].insertText(\"NF\u2265 5,0.\", Word.InsertLocation.replace);\n}\n\n// 5) Norma de recupera\u00e7\u00e3o ---------------------------------------------------\nconst norma = body.search(\n  \"M\u00e9dia aritm\u00e9tica da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recupera\u00e7\u00e3o.\",\n  { matchCase: true }\n);\nnorma.load(\"items\");\nawait context.sync();\nif (norma.items.length > 0) {\n  norma.items[0].insertText(\n    \"(NF+RC)/2 \u2265 5,0, onde RC \u00e9 uma prova de recupera\u00e7\u00e3o a ser aplicada.\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Apply the LOQ4206 course-plan updates:\n#  1. Activation date 2018 -> 2021\n#  2. Second \"Programa\" English summary paragraph expanded with full syllabus detail\n#  3. \"M\u00e9todo\" evaluation description replaced\n#  4. \"Crit\u00e9rio\" formula replaced\n#  5. \"Norma de recupera\u00e7\u00e3o\" text replaced\n\n$d = $word.ActiveDocument\n\n# 1) Ativa\u00e7\u00e3o date --------------------------------------------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"Ativa\u00e7\u00e3o: 01/01/2018\", $false, $false, $false, $false, $false, $true, 1, $false, \"Ativa\u00e7\u00e3o: 01/01/2021\", 1) | Out-Null\n\n# 2) Expand the second occurrence of the English \"Programa\" summary ------\n# (identical text also appears once under \"Programa resumido\" - that one\n# must stay untouched, so we walk the Paragraphs collection and update\n# only the *second* match.)\n$shortProgramaText = \"Introduction to Operational Research, Linear Programming, Simplex Method, Introduction to Graphs and Network Optimization, Case Study in Linear Programming, Introduction to Queue Theory.\"\n$fullProgramaText = \"1. Introduction to Operational Research 1.1. Concepts of Operational Research; 1.2. Modeling; 1.3. Structure of Mathematical Models; 1.4. Mathematical techniques in Operational Research; 1.2. Phases of a Study in Operational Research 2. Linear Programming 2.1. Definition 2.2. Formulation of Models 2.3. Graphic Resolution; 3. Simplex method 3.1. Development of the Simplex Method; 3.2. Simplex Method Procedure; 4. Introduction to Graphs and Network Optimization 4.1. Basic Concepts in Graph Theory 4.2. Maximum Flow Problems; 4.3. Minimum Path Problems 5. Case Studies in Linear Programming 5.1. Simple Transport Model 5.2. Model of Designation. 6. Introduction to Queuing Theory 6.1. Queuing Theory Concepts 6.2. Markovian Models\"\n\n$seen = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $shortProgramaText) {\n        $seen = $seen + 1\n        if ($seen -eq 2) {\n            $r = $p.Range\n            $r.MoveEnd(1, -1) | Out-Null\n            $r.Text = $fullProgramaText\n            break\n        }\n    }\n}\n\n# 3) M\u00e9todo description ----------------------------------------------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de laborat\u00f3rio, aulas de exerc\u00edcios.\", $false, $false, $false, $false, $false, $true, 1, $false, \"NF=A avalia\u00e7\u00e3o ser\u00e1 composta por provas, listas, projetos, semin\u00e1rios e outras formas que far\u00e3o a composi\u00e7\u00e3o das notas, sendo estipulada a m\u00e9dia final a somat\u00f3ria destas notas (N), com no m\u00ednimo duas avalia\u00e7\u00f5es, sendo: (N1+...+Nn)/n.\", 1) | Out-Null\n\n# 4) Crit\u00e9rio formula --------------------------------------------------------\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Execute(\"MF = (0,45*P1 + 0,45*P2 + 0,10*TRAB), onde P1 e P2 s\u00e3o provas e TRAB \u00e9 a nota m\u00e9dia de trabalhos e semin\u00e1rios.\", $false, $false, $false, $false, $false, $true, 1, $false, \"NF\u2265 5,0.\", 1) | Out-Null\n\n# 5) Norma de recupera\u00e7\u00e3o -----------------------------------------------------\n$rng4 = $d.Content\n$rng4.Find.ClearFormatting()\n$rng4.Find.Execute(\"M\u00e9dia aritm\u00e9tica da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recupera\u00e7\u00e3o.\", $false, $false, $false, $false, $false, $true, 1, $false, \"(NF+RC)/2 \u2265 5,0, onde RC \u00e9 uma prova de recupera\u00e7\u00e3o a ser aplicada.\", 1) | Out-Null\n"}
